$d = $word.ActiveDocument
$d.Content.Find.Execute("Recomened", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Recommend", 2)
